$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-reporting record was inserted as row 29 ("Rabanito" at
# "Región Metropolitana", $/paquete). This pushes every following record
# down by one row (old row 29 becomes 30, ..., old row 60 becomes 61).
$ws.Rows.Item(29).Insert()

# Fill in the newly inserted row 29 with its data.
$ws.Cells.Item(29, 1).Value = 10
$ws.Cells.Item(29, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(29, 3).Value = 'La Araucanía'
$ws.Cells.Item(29, 4).Value = 44664
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(29, 6).Value = 300000001
$ws.Cells.Item(29, 7).Value = 'Rabanito'
$ws.Cells.Item(29, 8).Value = 'Sin especificar'
$ws.Cells.Item(29, 9).Value = 'Primera'
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 1000
$ws.Cells.Item(29, 12).Value = 1000
$ws.Cells.Item(29, 13).Value = 1000
$ws.Cells.Item(29, 14).Value = '$/paquete'
$ws.Cells.Item(29, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 16).Value = 1000
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = 'Hortaliza'
